# Updated cryptos list -- applies the diff's cell-level changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / URL / percentage cells: safe to assign directly. ---
$ws.Range("D2").Value = "51.297.54"
$ws.Range("E2").Value = "  -1.91%  "
$ws.Range("D3").Value = "2.920.93"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  +2.01%  "
$ws.Range("E6").Value = "  -6.29%  "
$ws.Range("E7").Value = "  -4.87%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -6.86%  "
$ws.Range("E10").Value = "  -5.12%  "
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("E12").Value = "  -4.08%  "
$ws.Range("E13").Value = "  -5.79%  "
$ws.Range("D14").Value = "3.377.00"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("E15").Value = "  -5.25%  "
$ws.Range("D16").Value = "2.912.11"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("E17").Value = "  -3.50%  "
$ws.Range("D18").Value = "51.220.56"
$ws.Range("E18").Value = "  -1.78%  "
$ws.Range("E19").Value = "  -6.51%  "
$ws.Range("E20").Value = "  -3.84%  "
$ws.Range("E21").Value = "  -6.56%  "
$ws.Range("D22").Value = "0.0₃0947"
$ws.Range("E22").Value = "  -3.74%  "
$ws.Range("E23").Value = "  -3.52%  "
$ws.Range("E24").Value = "  -3.33%  "
$ws.Range("E25").Value = "  -4.42%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E26").Value = "  -5.02%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E28").Value = "  -3.78%  "
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E29").Value = "  -6.25%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("E30").Value = "  -3.57%  "
$ws.Range("B31").Value = "RenderToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("E32").Value = "  -5.50%  "
$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("E33").Value = "  -3.18%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("E34").Value = "  -6.53%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("E35").Value = "  -2.89%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E37").Value = "  -4.75%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("E38").Value = "  +2.39%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("E39").Value = "  -1.55%  "
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("E40").Value = "  -7.31%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("E41").Value = "  -6.54%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("E42").Value = "  -5.27%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E43").Value = "  -2.46%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("E44").Value = "  -1.52%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("E45").Value = "  -2.88%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.064.87"
$ws.Range("E46").Value = "  -3.09%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("E47").Value = "  -7.79%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("E48").Value = "  -8.46%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "3.203.33"
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("E50").Value = "  -6.69%  "
$ws.Range("E51").Value = "  -8.05%  "

# --- Price cells whose new text looks like a plain number (e.g. "364.72").
# Excel auto-converts such strings to Numbers on assignment, which would
# change the cell's stored type away from Text. Force Text by applying a
# "@" number format just for the assignment, then reset the cell style back
# to Normal so no extra formatting is left behind on the cell. ---
$numericTextCells = @{
    "D5" = "364.72"
    "D7" = "0.541"
    "D10" = "37.09"
    "D12" = "0.0836"
    "D13" = "18.47"
    "D15" = "7.37"
    "D17" = "0.957"
    "D20" = "7.29"
    "D23" = "68.25"
    "D24" = "260.51"
    "D26" = "0.175"
    "D27" = "1.00"
    "D28" = "25.99"
    "D29" = "7.29"
    "D30" = "0.104"
    "D31" = "6.17"
    "D32" = "9.95"
    "D33" = "2.14"
    "D34" = "35.22"
    "D35" = "50.73"
    "D36" = "1.00"
    "D37" = "0.0423"
    "D38" = "2.84"
    "D39" = "3.16"
    "D40" = "17.00"
    "D41" = "1.87"
    "D42" = "0.114"
    "D43" = "22.59"
    "D44" = "118.34"
    "D45" = "2.12"
    "D47" = "3.21"
    "D48" = "2.27"
    "D50" = "0.236"
    "D51" = "0.0318"
}
foreach ($addr in $numericTextCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextCells[$addr]
    $cell.Style = "Normal"
}
